$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update numeric values in the "Values" row (row 4) and "Defaults" row (row 3)
$ws.Range("D3").Value = 0.5
$ws.Range("C4").Value = 50000
$ws.Range("D4").Value = 0.5
$ws.Range("I4").Value = "30"
$ws.Range("J4").Value = 0.00005

# Update the selected cell in the sheet view
$ws.Range("D7").Select()
